$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new value in D7, reusing the same text as the header cell D1 ("Grupo_4")
$ws.Range("D7").Value = "Grupo_4"

# Update the active selection to D7, as reflected in the sheetView
$ws.Range("D7").Select()
